$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace "Freelance" income with "Bakery" income, and store the date
# as plain text (not an Excel date serial) to match the source data export.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("A2").Value = "Bakery"
$ws.Range("B2").Value = 20000
$ws.Range("C2").Value = "15/08/2025"
$ws.Range("C2").Style = "Normal"

# New rows 3-5: additional income entries, dates also stored as plain text.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("A3").Value = "MusicPlay"
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = "11/08/2025"
$ws.Range("C3").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("A4").Value = "Tractor"
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = "10/08/2025"
$ws.Range("C4").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("A5").Value = "Dukaan"
$ws.Range("B5").Value = 10000
$ws.Range("C5").Value = "04/07/2025"
$ws.Range("C5").Style = "Normal"
